# Doing Updates for Financials
# Insert a new "year" column (FY2018) in front of the existing data columns
# on the PBT sheet. This pushes the existing D:K columns (2017..2011) one
# column to the right (E:L) and populates the new column D with the FY2018
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D - this shifts D:K -> E:L (values + formatting)
# exactly like the Excel UI "Insert Sheet Columns" command.
$ws.Columns("D").Insert()

# The freshly inserted column D has no formatting yet; clone it from the
# column immediately to its right (the former column D, now E) so number
# formats / styles (date format, right-aligned numeric format, etc.) match.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# --- Income Statement (FY2018 column) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 32100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1300
$ws.Range("D18").Value = 30800
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 30800
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 30800
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 30800
$ws.Range("D27").Value = 30800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 30800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 30800

# --- Balance Sheet (FY2018 column) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 3500
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 500
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 4000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 2500
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 500
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (FY2018 column) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 30800
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 30700
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = -30800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -30800
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -100
